$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 96: 494. Target Sum ----
$ws.Range("A96").Value = "494. Target Sum"
$ws.Range("B96").Value = "Medium"
$ws.Range("C96").Value = "Dynamic Programming"
$ws.Range("D96").Value = "At each step (element in nums), there are 2 choices: add the element or subtract it. This forms a binary tree of decisions. We use a hashmap cache to store computations. The base case is when i == nums.length, and if sum equals the target, it is a valid way. The recursive step at each step should add the sum of both recursive calls."
$ws.Range("E96").Value = "https://leetcode.com/problems/target-sum/solutions/455024/dp-is-easy-5-steps-to-think-through-dp-questions/ "

# ---- Row 97: 190. Reverse Bits ----
$ws.Range("A97").Value = "190. Reverse Bits"
$ws.Range("B97").Value = "Easay"
$ws.Range("C97").Value = "Bit Manipulation"
$ws.Range("E97").Value = "https://leetcode.com/problems/reverse-bits/solutions/54738/sharing-my-2ms-java-solution-with-explanation/ "
$ws.Range("D97").Value = "Review bit manipulation. To get the current bit, we can AND with 1. To shift the 1 to see the bits to the left, we shift the 1 to the left with <<. To fill the output, we start by shifting 1 to the very left with << 31, then use OR with (n & 1) to replace the bits with the LSB of n, reversed. The crux is to understand how to work with least significant bits (LSB) as a pointer, left and right shift, Logic OR and Logic AND to get and replace LSBs."

# ---- Fill colors for Difficulty column (B), matching existing convention ----
# Medium = orange fill (same as B2)
$ws.Range("B96").Interior.Color = 49407
# Easy("Easay") = green fill (same as B7)
$ws.Range("B97").Interior.Color = 5287936

# ---- Hyperlinks for column E ----
$ws.Hyperlinks.Add($ws.Range("E96"), "https://leetcode.com/problems/target-sum/solutions/455024/dp-is-easy-5-steps-to-think-through-dp-questions/")
$ws.Range("E96").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E97"), "https://leetcode.com/problems/reverse-bits/solutions/54738/sharing-my-2ms-java-solution-with-explanation/")
$ws.Range("E97").Style = "Hyperlink"

# ---- Resize the table / autofilter to include the two new rows ----
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E97"))

# ---- Update selection to match the final view state ----
$ws.Range("E103").Select() | Out-Null
